$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.045608419067008
$ws.Range("D2").Value = 1.048080382523962
$ws.Range("E2").Value = 1.043079053674487
$ws.Range("F2").Value = 1.044447821620052
$ws.Range("I2").Value = 1.044609206762441
$ws.Range("J2").Value = 1.050667243542235
$ws.Range("K2").Value = 1.050841225952017
$ws.Range("L2").Value = 1.045853932001868
$ws.Range("M2").Value = 1.047218844371849

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.047079355999303
$ws.Range("D3").Value = 1.049255227248778
$ws.Range("E3").Value = 1.044372084595887
$ws.Range("F3").Value = 1.046546796370809
$ws.Range("I3").Value = 1.045159702223752
$ws.Range("J3").Value = 1.051783338939877
$ws.Range("K3").Value = 1.051827104596299
$ws.Range("L3").Value = 1.046956663991478
$ws.Range("M3").Value = 1.049125703094741

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.048028127108215
$ws.Range("D4").Value = 1.050012815318831
$ws.Range("E4").Value = 1.045205566982231
$ws.Range("F4").Value = 1.047901574026911
$ws.Range("I4").Value = 1.04551302914096
$ws.Range("J4").Value = 1.052502111585901
$ws.Range("K4").Value = 1.052461852011942
$ws.Range("L4").Value = 1.047666472505621
$ws.Range("M4").Value = 1.050355808828113

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.048426278985596
$ws.Range("D5").Value = 1.050330688430902
$ws.Range("E5").Value = 1.04555520894601
$ws.Range("F5").Value = 1.048470330578879
$ws.Range("I5").Value = 1.045660883694323
$ws.Range("J5").Value = 1.052803475951626
$ws.Range("K5").Value = 1.052727946741117
$ws.Range("L5").Value = 1.047963992580435
$ws.Range("M5").Value = 1.05087206533389

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.048493089071074
$ws.Range("D6").Value = 1.050384024812592
$ws.Range("E6").Value = 1.045613871467582
$ws.Range("F6").Value = 1.048565781495502
$ws.Range("I6").Value = 1.045685669168332
$ws.Range("J6").Value = 1.052854029284221
$ws.Range("K6").Value = 1.052772581306626
$ws.Range("L6").Value = 1.04801389600083
$ws.Range("M6").Value = 1.050958696051607

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.048033450013649
$ws.Range("D7").Value = 1.050017065167761
$ws.Range("E7").Value = 1.045210241864725
$ws.Range("F7").Value = 1.047909176864927
$ws.Range("I7").Value = 1.045515007461984
$ws.Range("J7").Value = 1.052506141592488
$ws.Range("K7").Value = 1.052465410531197
$ws.Range("L7").Value = 1.047670451442131
$ws.Range("M7").Value = 1.050362710505809

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.046106162565494
$ws.Range("D8").Value = 1.048477973780898
$ws.Range("E8").Value = 1.043516707668985
$ws.Range("F8").Value = 1.045157897029928
$ws.Range("I8").Value = 1.044795849028817
$ws.Range("J8").Value = 1.051045146544227
$ws.Range("K8").Value = 1.051175072365081
$ws.Range("L8").Value = 1.046227384526201
$ws.Range("M8").Value = 1.047864066734268

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.042686323883652
$ws.Range("D9").Value = 1.045745457720054
$ws.Range("E9").Value = 1.040507522666603
$ws.Range("F9").Value = 1.040282770532947
$ws.Range("I9").Value = 1.043506290413035
$ws.Range("J9").Value = 1.048444071644784
$ws.Range("K9").Value = 1.048876568890984
$ws.Range("L9").Value = 1.043655466099147
$ws.Range("M9").Value = 1.043431440366487

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.040389706708624
$ws.Range("D10").Value = 1.04390945621971
$ws.Range("E10").Value = 1.038483925979758
$ws.Range("F10").Value = 1.037013136474426
$ws.Range("I10").Value = 1.042631234397536
$ws.Range("J10").Value = 1.04669148720943
$ws.Range("K10").Value = 1.047327031238537
$ws.Range("L10").Value = 1.041920673191322
$ws.Range("M10").Value = 1.040455119132724

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.039391103837639
$ws.Range("D11").Value = 1.043110919184227
$ws.Range("E11").Value = 1.037603386859923
$ws.Range("F11").Value = 1.035592402469882
$ws.Range("I11").Value = 1.042248605766154
$ws.Range("J11").Value = 1.045928058986171
$ws.Range("K11").Value = 1.046651860020345
$ws.Range("L11").Value = 1.041164556553122
$ws.Range("M11").Value = 1.039161017715103

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.039019539173341
$ws.Range("D12").Value = 1.042813764316127
$ws.Range("E12").Value = 1.037275654665206
$ws.Range("F12").Value = 1.035063906518434
$ws.Range("I12").Value = 1.042105914145073
$ws.Range("J12").Value = 1.045643792029412
$ws.Range("K12").Value = 1.046400427866246
$ws.Range("L12").Value = 1.04088294634763
$ws.Range("M12").Value = 1.038679504740737

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.039099270222074
$ws.Range("D13").Value = 1.042877529784636
$ws.Range("E13").Value = 1.037345984452274
$ws.Range("F13").Value = 1.035177306065784
$ws.Range("I13").Value = 1.04213654772035
$ws.Range("J13").Value = 1.045704799941843
$ws.Range("K13").Value = 1.046454390211755
$ws.Range("L13").Value = 1.040943387054764
$ws.Range("M13").Value = 1.038782828688911

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.03936040331271
$ws.Range("D14").Value = 1.043086367406515
$ws.Range("E14").Value = 1.037576309972282
$ws.Range("F14").Value = 1.035548732736011
$ws.Range("I14").Value = 1.042236822424939
$ws.Range("J14").Value = 1.045904575656244
$ws.Range("K14").Value = 1.046631089750366
$ws.Range("L14").Value = 1.041141294037581
$ws.Range("M14").Value = 1.039121232697872

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.03952121088313
$ws.Range("D15").Value = 1.043214966907616
$ws.Range("E15").Value = 1.037718133163988
$ws.Range("F15").Value = 1.035777477850961
$ws.Range("I15").Value = 1.042298529730846
$ws.Range("J15").Value = 1.0460275715007
$ws.Range("K15").Value = 1.046739874571711
$ws.Range("L15").Value = 1.041263130659058
$ws.Range("M15").Value = 1.03932962437023

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.040455891908479
$ws.Range("D16").Value = 1.043962377020215
$ws.Range("E16").Value = 1.038542272542003
$ws.Range("F16").Value = 1.037107319114064
$ws.Range("I16").Value = 1.042656549162241
$ws.Range("J16").Value = 1.046742056550456
$ws.Range("K16").Value = 1.047371750483498
$ws.Range("L16").Value = 1.041970749007943
$ws.Range("M16").Value = 1.040540889892405

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.04104107018985
$ws.Range("D17").Value = 1.04443025289377
$ws.Range("E17").Value = 1.039058070460103
$ws.Range("F17").Value = 1.037940146444069
$ws.Range("I17").Value = 1.042880123548029
$ws.Range("J17").Value = 1.04718900785037
$ws.Range("K17").Value = 1.047766974499816
$ws.Range("L17").Value = 1.042413287871234
$ws.Range("M17").Value = 1.041299239228586

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.041381995165396
$ws.Range("D18").Value = 1.044702816552196
$ws.Range("E18").Value = 1.039358511879342
$ws.Range("F18").Value = 1.038425443706656
$ws.Range("I18").Value = 1.043010171900552
$ws.Range("J18").Value = 1.047449269056846
$ws.Range("K18").Value = 1.047997096401105
$ws.Range("L18").Value = 1.042670937541209
$ws.Range("M18").Value = 1.041741058524144

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.041498174448037
$ws.Range("D19").Value = 1.044795696316228
$ws.Range("E19").Value = 1.03946088469039
$ws.Range("F19").Value = 1.038590837458625
$ws.Range("I19").Value = 1.043054454375422
$ws.Range("J19").Value = 1.04753793749855
$ws.Range("K19").Value = 1.048075493554177
$ws.Range("L19").Value = 1.04275870910367
$ws.Range("M19").Value = 1.041891621205673

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.040978327542357
$ws.Range("D20").Value = 1.044380089534075
$ws.Range("E20").Value = 1.039002773224792
$ws.Range("F20").Value = 1.037850841466937
$ws.Range("I20").Value = 1.042856173294769
$ws.Range("J20").Value = 1.0471410995926
$ws.Range("K20").Value = 1.047724612732406
$ws.Range("L20").Value = 1.042365856940996
$ws.Range("M20").Value = 1.041217928765525

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.039283523891497
$ws.Range("D21").Value = 1.043024885004103
$ws.Range("E21").Value = 1.037508503196551
$ws.Range("F21").Value = 1.03543937836913
$ws.Range("I21").Value = 1.042207309712373
$ws.Range("J21").Value = 1.045845765994135
$ws.Range("K21").Value = 1.046579073991722
$ws.Range("L21").Value = 1.041083036303016
$ws.Range("M21").Value = 1.039021604154994

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.038214228281512
$ws.Range("D22").Value = 1.042169668734791
$ws.Range("E22").Value = 1.036565166579316
$ws.Range("F22").Value = 1.03391871526484
$ws.Range("I22").Value = 1.04179606345203
$ws.Range("J22").Value = 1.045027306806334
$ws.Range("K22").Value = 1.045855099758329
$ws.Range("L22").Value = 1.040272102514071
$ws.Range("M22").Value = 1.037635895793614

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.038781438468125
$ws.Range("D23").Value = 1.042623337042272
$ws.Range("E23").Value = 1.037065614723797
$ws.Range("F23").Value = 1.034725281276493
$ws.Range("I23").Value = 1.042014386099494
$ws.Range("J23").Value = 1.045461573864113
$ws.Range("K23").Value = 1.046239249269266
$ws.Range("L23").Value = 1.040702412686787
$ws.Range("M23").Value = 1.038370948529391

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.041006679478637
$ws.Range("D24").Value = 1.044402757248873
$ws.Range("E24").Value = 1.039027760946687
$ws.Range("F24").Value = 1.037891196014701
$ws.Range("I24").Value = 1.042866996491891
$ws.Range("J24").Value = 1.047162748625098
$ws.Range("K24").Value = 1.047743755444758
$ws.Range("L24").Value = 1.042387290406636
$ws.Range("M24").Value = 1.041254671048467

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.043573324109701
$ws.Range("D25").Value = 1.046454360300126
$ws.Range("E25").Value = 1.041288497790537
$ws.Range("F25").Value = 1.041546450985279
$ws.Range("I25").Value = 1.043842351847972
$ws.Range("J25").Value = 1.049119731292235
$ws.Range("K25").Value = 1.049473778155196
$ws.Range("L25").Value = 1.044323877600715
$ws.Range("M25").Value = 1.044581029646666
